$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reproduce the exact shared-string insertion order used by the original author.
$ws.Range("B5").Value = "https://github.com/anshupandey/Computer-Vision/blob/master/Transfer%20Learning%20for%20CNNs/Transfer_Learning.pptx"
$ws.Range("A5").Value = "Transfer Learning for computer vision"
$ws.Range("A7").Value = "word2vec"
$ws.Range("A8").Value = "transfer learning for NLP"
$ws.Range("A9").Value = "Attention based network, transformers"
$ws.Range("A6").Value = "RNN & LSTM"

$ws.Range("B6").Value = "https://github.com/anshupandey/Natural_language_Processing/blob/master/NLP%20-%20Deep%20Learning%20-%20RNN%20%26%20LSTM/RNN.pptx"
$ws.Range("B7").Value = "https://github.com/anshupandey/Natural_language_Processing/blob/master/NLP%20-%20Deep%20Learning%20-%20RNN%20%26%20LSTM/Word2vec.pptx"
$ws.Range("B8").Value = "https://github.com/anshupandey/Natural_language_Processing/blob/master/NLP%20-%20Transfer%20Learning/Transformer%20Model.pptx"
$ws.Range("B9").Value = "https://github.com/anshupandey/Natural_language_Processing/blob/master/NLP%20-%20Transfer%20Learning/Transfer%20Learning%20in%20NLP.pptx"

# Apply hyperlinks for rows 6, 8 and 9 (rows 5 and 7 stay as plain text URLs).
$ws.Hyperlinks.Add($ws.Range("B6"), "https://github.com/anshupandey/Natural_language_Processing/blob/master/NLP%20-%20Deep%20Learning%20-%20RNN%20%26%20LSTM/RNN.pptx") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B8"), "https://github.com/anshupandey/Natural_language_Processing/blob/master/NLP%20-%20Transfer%20Learning/Transformer%20Model.pptx") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B9"), "https://github.com/anshupandey/Natural_language_Processing/blob/master/NLP%20-%20Transfer%20Learning/Transfer%20Learning%20in%20NLP.pptx") | Out-Null

# Match final selection state from the target workbook.
$ws.Range("A9").Select() | Out-Null

Write-Host "Done applying PPT link updates"
